# Updated cryptos list on Wed May 24 07:09:27 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto table.
# Row 44/45 additionally swap (Decentraland <-> PaxDollar changed rank order),
# so their Coin name (B) and Link (C) cells are updated too.
#
# Note: values that are purely numeric-looking (e.g. "1.006") are written
# with a leading apostrophe so Excel keeps them as literal text instead of
# auto-converting to a number, matching the workbook's inlineStr cell type.
# Values that already contain stray characters Excel can't parse as a
# number (e.g. "26.866.11", two dots) are left unprefixed since Excel
# stores those as text on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.866.11'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '1.832.40'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''310.37'
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("D6").Value = '''1.005'
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("D8").Value = '''0.3679'
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '''0.07166'
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("D10").Value = '''0.8772'
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D11").Value = '''0.07867'
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Value = '''19.60'
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").Value = '1.853.40'
$ws.Range("E13").Value = '  -1.43%  '
$ws.Range("D14").Value = '''5.333'
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").Value = '''6.385'
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("D16").Value = '''86.89'
$ws.Range("E16").Value = '  -5.89%  '
$ws.Range("D17").Value = '''1.006'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '''0.000008721'
$ws.Range("E18").Value = '  -1.87%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '26.885.81'
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("D21").Value = '''14.45'
$ws.Range("E21").Value = '  -2.80%  '
$ws.Range("D22").Value = '''4.996'
$ws.Range("E22").Value = '  -2.99%  '
$ws.Range("D23").Value = '''10.44'
$ws.Range("E23").Value = '  -1.03%  '
$ws.Range("D24").Value = '''1.981'
$ws.Range("E24").Value = '  +4.47%  '
$ws.Range("D25").Value = '''150.70'
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("D26").Value = '''18.23'
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("D27").Value = '''1.970'
$ws.Range("E27").Value = '  -5.05%  '
$ws.Range("D28").Value = '''113.50'
$ws.Range("E28").Value = '  -2.61%  '
$ws.Range("D29").Value = '''4.931'
$ws.Range("E29").Value = '  -4.00%  '
$ws.Range("D30").Value = '''0.08816'
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").Value = '''3.126'
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("D32").Value = '''0.7552'
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("D33").Value = '''4.460'
$ws.Range("D34").Value = '''1.129'
$ws.Range("D35").Value = '''2.579'
$ws.Range("E35").Value = '  -2.88%  '
$ws.Range("D36").Value = '''1.087'
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").Value = '''0.01933'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").Value = '''2.929'
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("D39").Value = '''0.05124'
$ws.Range("E39").Value = '  -2.49%  '
$ws.Range("D40").Value = '''6.896'
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("D41").Value = '''0.4975'
$ws.Range("E41").Value = '  -3.88%  '
$ws.Range("D42").Value = '''0.1595'
$ws.Range("E42").Value = '  -3.06%  '
$ws.Range("D43").Value = '''8.335'
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '''0.4670'
$ws.Range("E44").Value = '  -3.75%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''1.005'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").Value = '''10.09'
$ws.Range("E46").Value = '  -2.68%  '
$ws.Range("D47").Value = '''102.30'
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("D48").Value = '''1.611'
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("D49").Value = '''0.06093'
$ws.Range("E49").Value = '  -2.48%  '
$ws.Range("D50").Value = '''64.41'
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = '''36.32'
$ws.Range("E51").Value = '  -2.35%  '
